$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data table: Row,F(serial or N),G(int or N),H(serial or N),I(int or N)
# Represents the recomputed bond coupon schedule after advancing "today"
# from 2023-10-19 to 2023-10-20 (one day later).
$data = @"
2,N,N,45340,121
3,45121,98,45305,86
4,45039,180,45222,3
5,45039,180,45222,3
6,45129,90,45221,2
7,45162,57,45254,35
8,45156,63,45248,29
9,45141,78,45233,14
10,45062,157,45246,27
11,N,N,45286,67
12,45195,24,45377,158
13,45036,183,45402,183
14,45036,183,45402,183
15,45106,113,45289,70
16,45186,33,45278,59
17,45186,33,45277,58
18,45132,87,45224,5
19,45137,82,45229,10
20,45183,36,45274,55
21,45187,32,45278,59
22,45187,32,45278,59
23,N,N,45245,26
24,N,N,45245,26
25,45204,15,N,N
26,45160,59,45344,125
27,45118,101,45302,83
28,45085,134,45268,49
29,45183,36,45274,55
30,N,N,45224,5
31,45153,66,45245,26
32,45155,64,45247,28
33,N,N,45255,36
34,45176,43,45267,48
35,45130,89,45222,3
36,45207,12,45299,80
37,45189,30,45280,61
38,45107,112,45290,71
39,45063,156,45247,28
40,45177,42,45359,140
41,45107,112,45290,71
42,N,N,45265,46
43,45098,121,45281,62
44,45099,120,45282,63
45,45106,113,45289,70
46,45106,113,45289,70
47,45120,99,45304,85
48,45118,101,45302,83
49,45152,67,45336,117
50,45042,177,45225,6
51,45122,97,45306,87
52,45139,80,45323,104
53,45099,120,45282,63
54,45196,23,45287,68
55,45044,175,45227,8
56,45151,68,45243,24
57,45196,23,45287,68
58,45069,150,45253,34
59,45191,28,45282,63
60,45205,14,N,N
61,45127,92,45313,94
62,45058,161,45242,23
63,45194,25,45376,157
64,45208,11,N,N
65,45208,11,N,N
66,45089,130,45272,53
67,45178,41,45269,50
68,45196,23,45287,68
69,45149,70,45241,22
70,45205,14,45296,77
71,45205,14,N,N
72,45207,12,45299,80
73,45189,30,45280,61
74,45038,181,45221,2
75,45117,102,45301,82
76,45160,59,45344,125
77,45186,33,45368,149
78,45043,176,45226,7
79,45079,140,45262,43
80,N,N,45254,35
81,45104,115,45287,68
82,N,N,45288,69
83,45175,44,45357,138
84,45054,165,45238,19
85,45056,163,45240,21
86,45085,134,45268,49
87,45194,25,45560,341
88,45103,116,45286,67
89,45103,116,45286,67
90,45192,27,45374,155
91,45100,119,45287,68
92,45187,32,45369,150
93,45104,115,45287,68
94,45104,115,45287,68
95,45049,170,45233,14
96,45122,97,45306,87
97,45141,78,45325,106
98,45194,25,45376,157
99,45093,126,45276,57
100,45123,96,45307,88
101,45198,21,45380,161
102,45093,126,45276,57
103,45122,97,45306,87
104,45122,97,45306,87
105,45050,169,45234,15
106,45085,134,45268,49
107,45185,34,45276,57
108,45185,34,45276,57
109,45185,34,45276,57
110,45148,71,45240,21
111,45211,8,45303,84
112,45172,47,45263,44
113,45206,13,45389,170
114,45203,16,45295,76
115,45218,1,45310,91
116,45161,58,45253,34
117,45115,104,45299,80
118,45094,125,45277,58
119,45104,115,45287,68
120,45104,115,45287,68
121,45104,115,45287,68
122,45104,115,45287,68
123,45199,20,45291,72
124,45197,22,45379,160
125,45213,6,45305,86
126,45134,85,45226,7
127,45212,7,45304,85
128,45218,1,45310,91
129,45158,61,45250,31
130,45171,48,45262,43
131,45203,16,45295,76
132,45144,75,45236,17
133,45156,63,45248,29
134,45146,73,45238,19
135,45140,79,45232,13
136,45196,23,45287,68
137,45197,22,45288,69
138,45179,40,45270,51
139,45197,22,45288,69
140,45182,37,45273,54
141,45208,11,45300,81
142,45150,69,45242,23
143,N,N,45238,19
144,45171,48,45262,43
145,45194,25,45285,66
146,45211,8,N,N
147,45165,54,45257,38
148,45166,53,45258,39
149,45185,34,45276,57
150,45218,1,45401,182
151,45174,45,45356,137
152,45089,130,45272,53
153,45124,95,45308,89
154,45124,95,45308,89
155,45209,10,45392,173
156,45209,10,45392,173
157,45209,10,45392,173
158,45209,10,45392,173
159,45156,63,45247,28
160,45153,66,45245,26
161,45135,84,45319,100
162,45091,128,45274,55
163,45091,128,45274,55
164,45138,81,45322,103
165,45138,81,45322,103
166,45085,134,45268,49
167,45085,134,45268,49
168,N,N,45326,107
169,45173,46,45355,136
170,45056,163,45240,21
171,N,N,45238,19
172,45141,78,45325,106
173,45037,182,45220,1
174,45171,48,45262,43
175,45197,22,45288,69
176,45162,57,45254,35
177,45061,158,45245,26
178,45212,7,45304,85
179,45191,28,45282,63
180,45088,131,45271,52
181,45192,27,45283,64
182,45131,88,45223,4
183,45136,83,45228,9
184,45150,69,45242,23
185,45210,9,N,N
186,45129,90,45221,2
187,45189,30,45280,61
188,45193,26,45284,65
189,45165,54,45257,38
190,45202,17,45385,166
191,45081,138,45264,45
192,45081,138,45264,45
193,N,N,45319,100
194,N,N,45319,100
195,45046,173,45230,11
196,45046,173,45230,11
197,45214,5,45397,178
198,45214,5,45397,178
199,45214,5,45397,178
200,45214,5,45397,178
201,45067,152,45251,32
202,45067,152,45251,32
203,45067,152,45251,32
204,45067,152,45251,32
205,45178,41,45269,50
206,N,N,45252,33
207,45087,132,45270,51
208,45099,120,45282,63
209,45175,44,45356,137
210,45175,44,45356,137
211,45166,53,45350,131
212,45166,53,45350,131
213,45090,129,45273,54
214,45106,113,45289,70
215,45134,85,45226,7
216,45132,87,45224,5
217,45150,69,45242,23
218,45187,32,45278,59
219,45182,37,45273,54
220,45196,23,45287,68
221,45200,19,45292,73
222,45135,84,45227,8
223,45148,71,45240,21
224,45192,27,45283,64
225,45136,83,45320,101
226,45136,83,45320,101
227,45094,125,45277,58
228,45081,138,45264,45
229,45099,120,45282,63
230,45201,18,45384,165
231,45214,5,45397,178
232,N,N,45294,75
233,45079,140,45262,43
234,45199,20,45382,163
235,45204,15,45387,168
236,45204,15,45387,168
237,45080,139,45263,44
238,45046,173,45229,10
239,45204,15,45296,77
240,45204,15,45296,77
241,45199,20,45382,163
242,45070,149,45254,35
243,45168,51,45260,41
244,45210,9,45393,174
245,45077,142,45260,41
246,45100,119,45283,64
247,N,N,45281,62
248,45044,175,45227,8
249,45101,118,45283,64
250,45134,85,45226,7
251,45158,61,45250,31
252,45200,19,45292,73
253,N,N,45317,98
254,45039,180,45222,3
255,45129,90,45221,2
256,45157,62,45341,122
257,45063,156,45247,28
258,N,N,45254,35
259,N,N,45486,267
260,45073,146,45439,220
261,N,N,45323,104
262,N,N,45303,84
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line.Length -eq 0) { continue }
    $parts = $line -split ","
    $row = [int]$parts[0]
    $fVal = $parts[1]
    $gVal = $parts[2]
    $hVal = $parts[3]
    $iVal = $parts[4]

    if ($fVal -ne "N") {
        $ws.Cells.Item($row, 6).Value2 = [double]$fVal
    }
    if ($gVal -ne "N") {
        $ws.Cells.Item($row, 7).Value2 = [double]$gVal
    }
    if ($hVal -ne "N") {
        $ws.Cells.Item($row, 8).Value2 = [double]$hVal
    }
    if ($iVal -ne "N") {
        $ws.Cells.Item($row, 9).Value2 = [double]$iVal
    }
}

Write-Host "Updated $($lines.Count) rows"
